# Update header row (row 1): MobileNo, OTP, Population, Value
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "MobileNo"
$ws.Range("B1").Value = "OTP"
$ws.Range("C1").Value = "Population"
$ws.Range("D1").Value = "Value "

# Update data row (row 2) with the new sample mobile number / OTP values
$ws.Range("A2").Value = 510000017
$ws.Range("B2").Value = 1212

# Move the active selection to B2
$ws.Range("B2").Select() | Out-Null
